$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = [double]"25.37000000000053"
$ws.Range("H2").Value = [double]"5.871149605374271e-07"
$ws.Range("I2").Value = [double]"5.871149605374271e-07"
$ws.Range("L2").Value = [double]"44.39729629181707"
$ws.Range("M2").Value = "[27.591528885605015, 61.203063698029126]"
$ws.Range("N2").Value = [double]"3.139109162830778e-06"
$ws.Range("O2").Value = [double]"3.139109162830778e-06"
$ws.Range("P2").Value = [double]"1.37739497724958"
$ws.Range("Q2").Value = "[0.9371317425122703, 1.8176582119868891]"
$ws.Range("R2").Value = [double]"1.113625549553632e-07"
$ws.Range("S2").Value = [double]"1.113625549553632e-07"
$ws.Range("T2").Value = [double]"52.51528898667804"
$ws.Range("U2").Value = "[42.338380291502446, 62.69219768185364]"
$ws.Range("V2").Value = [double]"1.532107773982716e-13"
$ws.Range("W2").Value = [double]"1.532107773982716e-13"
$ws.Range("X2").Value = [double]"19.80840840840882"
$ws.Range("Y2").Value = [double]"18.0307307307311"
$ws.Range("Z2").Value = [double]"21.58608608608654"

# Row 3
$ws.Range("F3").Value = [double]"25.37000000000053"
$ws.Range("H3").Value = [double]"2.6469689312969e-05"
$ws.Range("I3").Value = [double]"2.6469689312969e-05"
$ws.Range("L3").Value = [double]"40.41935850221949"
$ws.Range("M3").Value = "[19.538696347988065, 61.30002065645091]"
$ws.Range("N3").Value = [double]"0.0003191457541142917"
$ws.Range("O3").Value = [double]"0.0003191457541142917"
$ws.Range("P3").Value = [double]"1.037763339023655"
$ws.Range("Q3").Value = "[0.5220264069028095, 1.5535002711445012]"
$ws.Range("R3").Value = [double]"0.0001978673728328939"
$ws.Range("S3").Value = [double]"0.0001978673728328939"
$ws.Range("T3").Value = [double]"57.73972927906161"
$ws.Range("U3").Value = "[46.60717676366141, 68.8722817944618]"
$ws.Range("V3").Value = [double]"1.301181384860683e-13"
$ws.Range("W3").Value = [double]"1.301181384860683e-13"
$ws.Range("X3").Value = [double]"21.1797597597602"
$ws.Range("Y3").Value = [double]"19.09733733733774"
$ws.Range("Z3").Value = [double]"23.26218218218266"

# Row 4
$ws.Range("F4").Value = [double]"25.37000000000053"
$ws.Range("H4").Value = [double]"5.554633876192483e-07"
$ws.Range("I4").Value = [double]"5.554633876192483e-07"
$ws.Range("L4").Value = [double]"48.64246272076683"
$ws.Range("M4").Value = "[28.370373870296508, 68.91455157123715]"
$ws.Range("N4").Value = [double]"1.60176490104913e-05"
$ws.Range("O4").Value = [double]"1.60176490104913e-05"
$ws.Range("P4").Value = [double]"0.748447499053424"
$ws.Range("Q4").Value = "[0.34592111300788453, 1.1509738850989635]"
$ws.Range("R4").Value = [double]"0.0005106187349821489"
$ws.Range("S4").Value = [double]"0.0005106187349821489"
$ws.Range("T4").Value = [double]"54.01426003601598"
$ws.Range("U4").Value = "[43.26609382030155, 64.76242625173042]"
$ws.Range("V4").Value = [double]"3.554934124849751e-13"
$ws.Range("W4").Value = [double]"3.554934124849751e-13"
$ws.Range("X4").Value = [double]"22.34794794794841"
$ws.Range("Y4").Value = [double]"20.72264264264307"
$ws.Range("Z4").Value = [double]"23.97325325325375"

# Row 5
$ws.Range("B5").Value = 0
$ws.Range("F5").Value = [double]"25.37000000000053"
$ws.Range("H5").Value = [double]"1.912057046027371e-08"
$ws.Range("I5").Value = [double]"1.912057046027371e-08"
$ws.Range("L5").Value = [double]"49.25107858507421"
$ws.Range("M5").Value = "[32.89766583970413, 65.60449133044429]"
$ws.Range("N5").Value = [double]"2.49320043277379e-07"
$ws.Range("O5").Value = [double]"2.49320043277379e-07"
$ws.Range("P5").Value = [double]"0.3207632138800394"
$ws.Range("Q5").Value = "[-0.044026323473730145, 0.6855527512338089]"
$ws.Range("R5").Value = [double]"0.08332996393067682"
$ws.Range("S5").Value = [double]"0.08332996393067682"
$ws.Range("T5").Value = [double]"55.12351260380513"
$ws.Range("U5").Value = "[45.864532043176496, 64.38249316443377]"
$ws.Range("V5").Value = [double]"1.332267629550188e-15"
$ws.Range("W5").Value = [double]"1.332267629550188e-15"
$ws.Range("X5").Value = [double]"24.07483483483533"
$ws.Range("Y5").Value = [double]"22.60190190190237"
$ws.Range("Z5").Value = [double]"25.5477677677683"

# Row 6
$ws.Range("F6").Value = [double]"25.37000000000053"
$ws.Range("H6").Value = [double]"2.926720643614544e-09"
$ws.Range("I6").Value = [double]"2.926720643614544e-09"
$ws.Range("L6").Value = [double]"58.58175429540787"
$ws.Range("M6").Value = "[42.03783567156468, 75.12567291925106]"
$ws.Range("N6").Value = [double]"6.49720566237022e-09"
$ws.Range("O6").Value = [double]"6.49720566237022e-09"
$ws.Range("P6").Value = [double]"0.05660527303765406"
$ws.Range("Q6").Value = "[-0.25786846606042335, 0.37107901213573147]"
$ws.Range("R6").Value = [double]"0.7186467968442951"
$ws.Range("S6").Value = [double]"0.7186467968442951"
$ws.Range("T6").Value = [double]"53.17165459019312"
$ws.Range("U6").Value = "[42.92297534828192, 63.42033383210432]"
$ws.Range("V6").Value = [double]"1.287858708565182e-13"
$ws.Range("W6").Value = [double]"1.287858708565182e-13"
$ws.Range("X6").Value = [double]"25.14144144144196"
$ws.Range("Y6").Value = [double]"23.87167167167217"
$ws.Range("Z6").Value = [double]"26.41121121121176"

# Row 7
$ws.Range("F7").Value = [double]"25.37000000000053"
$ws.Range("H7").Value = [double]"1.029661625961964e-07"
$ws.Range("I7").Value = [double]"1.029661625961964e-07"
$ws.Range("L7").Value = [double]"57.18603574441786"
$ws.Range("M7").Value = "[36.37078608417691, 78.0012854046588]"
$ws.Range("N7").Value = [double]"1.530162107910016e-06"
$ws.Range("O7").Value = [double]"1.530162107910016e-06"
$ws.Range("P7").Value = [double]"-0.3773684869176925"
$ws.Range("Q7").Value = "[-0.779894872963232, 0.02515789912784694]"
$ws.Range("R7").Value = [double]"0.06545332432287432"
$ws.Range("S7").Value = [double]"0.06545332432287432"
$ws.Range("T7").Value = [double]"60.7816870303252"
$ws.Range("U7").Value = "[49.01323117999188, 72.55014288065853]"
$ws.Range("V7").Value = [double]"1.48991929904696e-13"
$ws.Range("W7").Value = [double]"1.48991929904696e-13"
$ws.Range("X7").Value = [double]"1.523723723723755"
$ws.Range("Y7").Value = [double]"-0.1015815815815857"
$ws.Range("Z7").Value = [double]"3.149029029029096"

# Row 8
$ws.Range("F8").Value = [double]"25.49000000000055"
$ws.Range("H8").Value = [double]"1.686988083227092e-06"
$ws.Range("I8").Value = [double]"1.686988083227092e-06"
$ws.Range("L8").Value = [double]"43.13366381181657"
$ws.Range("M8").Value = "[23.248006637670265, 63.01932098596287]"
$ws.Range("N8").Value = [double]"7.267843498759419e-05"
$ws.Range("O8").Value = [double]"7.267843498759419e-05"
$ws.Range("P8").Value = [double]"-0.8050527720910781"
$ws.Range("Q8").Value = "[-1.2453160068283866, -0.3647895373537695]"
$ws.Range("R8").Value = [double]"0.0006158372908429133"
$ws.Range("S8").Value = [double]"0.0006158372908429133"
$ws.Range("T8").Value = [double]"42.73246868873241"
$ws.Range("U8").Value = "[32.3464544816149, 53.11848289584992]"
$ws.Range("V8").Value = [double]"1.322948417481484e-10"
$ws.Range("W8").Value = [double]"1.322948417481484e-10"
$ws.Range("X8").Value = [double]"3.265985985986056"
$ws.Range("Y8").Value = [double]"1.479899899899931"
$ws.Range("Z8").Value = [double]"5.052072072072181"

# Row 9
$ws.Range("F9").Value = [double]"25.49000000000055"
$ws.Range("H9").Value = [double]"1.017038342432386e-09"
$ws.Range("I9").Value = [double]"1.017038342432386e-09"
$ws.Range("L9").Value = [double]"59.30855733594404"
$ws.Range("M9").Value = "[41.10025647031915, 77.51685820156892]"
$ws.Range("N9").Value = [double]"4.584112445193966e-08"
$ws.Range("O9").Value = [double]"4.584112445193966e-08"
$ws.Range("P9").Value = [double]"-1.308210754648002"
$ws.Range("Q9").Value = "[-1.6604213424378482, -0.9560001668581553]"
$ws.Range("R9").Value = [double]"1.983203157251978e-09"
$ws.Range("S9").Value = [double]"1.983203157251978e-09"
$ws.Range("T9").Value = [double]"56.28042450220811"
$ws.Range("U9").Value = "[45.698512454139184, 66.86233655027704]"
$ws.Range("V9").Value = [double]"5.773159728050814e-14"
$ws.Range("W9").Value = [double]"5.773159728050814e-14"
$ws.Range("X9").Value = [double]"5.30722722722734"
$ws.Range("Y9").Value = [double]"3.878358358358442"
$ws.Range("Z9").Value = [double]"6.736096096096238"

# Row 10
$ws.Range("F10").Value = [double]"25.49000000000055"
$ws.Range("H10").Value = [double]"4.200755965300118e-09"
$ws.Range("I10").Value = [double]"4.200755965300118e-09"
$ws.Range("L10").Value = [double]"64.69030986093142"
$ws.Range("M10").Value = "[46.27458515324497, 83.10603456861787]"
$ws.Range("N10").Value = [double]"7.886195696471532e-09"
$ws.Range("O10").Value = [double]"7.886195696471532e-09"
$ws.Range("P10").Value = [double]"-1.622684493746079"
$ws.Range("Q10").Value = "[-1.9748950815359265, -1.2704739059562318]"
$ws.Range("R10").Value = [double]"5.133227176656874e-12"
$ws.Range("S10").Value = [double]"5.133227176656874e-12"
$ws.Range("T10").Value = [double]"74.04764363696621"
$ws.Range("U10").Value = "[61.92302939603678, 86.17225787789565]"
$ws.Range("V10").Value = [double]"6.661338147750939e-16"
$ws.Range("W10").Value = [double]"6.661338147750939e-16"
$ws.Range("X10").Value = [double]"6.583003003003142"
$ws.Range("Y10").Value = [double]"5.15413413413424"
$ws.Range("Z10").Value = [double]"8.011871871872044"

# Row 11
$ws.Range("F11").Value = [double]"25.49000000000055"
$ws.Range("H11").Value = [double]"0.0002893302362854122"
$ws.Range("I11").Value = [double]"0.0002893302362854122"
$ws.Range("L11").Value = [double]"39.54095147450269"
$ws.Range("M11").Value = "[16.355001109789242, 62.72690183921614]"
$ws.Range("N11").Value = [double]"0.001284811832241939"
$ws.Range("O11").Value = [double]"0.001284811832241939"
$ws.Range("P11").Value = [double]"-1.798789787641002"
$ws.Range("Q11").Value = "[-2.4654741145289267, -1.132105460753078]"
$ws.Range("R11").Value = [double]"2.140280573659226e-06"
$ws.Range("S11").Value = [double]"2.140280573659226e-06"
$ws.Range("T11").Value = [double]"53.06103447520897"
$ws.Range("U11").Value = "[40.10469916435986, 66.01736978605808]"
$ws.Range("V11").Value = [double]"1.503084323672965e-10"
$ws.Range("W11").Value = [double]"1.503084323672965e-10"
$ws.Range("X11").Value = [double]"7.297437437437594"
$ws.Range("Y11").Value = [double]"4.59279279279289"
$ws.Range("Z11").Value = [double]"10.0020820820823"

# Row 12
$ws.Range("F12").Value = [double]"25.49000000000055"
$ws.Range("H12").Value = [double]"4.338682191296073e-08"
$ws.Range("I12").Value = [double]"4.338682191296073e-08"
$ws.Range("L12").Value = [double]"53.79532184428711"
$ws.Range("M12").Value = "[33.223361861880676, 74.36728182669354]"
$ws.Range("N12").Value = [double]"3.765006227007817e-06"
$ws.Range("O12").Value = [double]"3.765006227007817e-06"
$ws.Range("P12").Value = [double]"-2.364842518017542"
$ws.Range("Q12").Value = "[-2.742211004935235, -1.987474031099849]"
$ws.Range("R12").Value = [double]"2.220446049250313e-16"
$ws.Range("S12").Value = [double]"2.220446049250313e-16"
$ws.Range("T12").Value = [double]"64.1715450724639"
$ws.Range("U12").Value = "[53.56496726203143, 74.77812288289637]"
$ws.Range("V12").Value = [double]"8.881784197001252e-16"
$ws.Range("W12").Value = [double]"8.881784197001252e-16"
$ws.Range("X12").Value = [double]"9.59383383383404"
$ws.Range("Y12").Value = [double]"8.062902902903076"
$ws.Range("Z12").Value = [double]"11.124764764765"

# Row 13
$ws.Range("F13").Value = [double]"25.49000000000055"
$ws.Range("H13").Value = [double]"3.978928297954099e-12"
$ws.Range("I13").Value = [double]"3.978928297954099e-12"
$ws.Range("L13").Value = [double]"66.98311993614496"
$ws.Range("M13").Value = "[50.28489231118212, 83.6813475611078]"
$ws.Range("N13").Value = [double]"2.642912555472776e-10"
$ws.Range("O13").Value = [double]"2.642912555472776e-10"
$ws.Range("P13").Value = [double]"-2.779947853627004"
$ws.Range("Q13").Value = "[-3.0566847440333125, -2.503210963220696]"
$ws.Range("T13").Value = [double]"58.93211062015256"
$ws.Range("U13").Value = "[49.65943072014146, 68.20479052016366]"
$ws.Range("V13").Value = [double]"0"
$ws.Range("W13").Value = [double]"0"
$ws.Range("X13").Value = [double]"11.2778578578581"
$ws.Range("Y13").Value = [double]"10.15517517517539"
$ws.Range("Z13").Value = [double]"12.40054054054081"

# Row 14
$ws.Range("F14").Value = [double]"25.49000000000055"
$ws.Range("H14").Value = [double]"1.494009693736587e-09"
$ws.Range("I14").Value = [double]"1.494009693736587e-09"
$ws.Range("L14").Value = [double]"61.01472230287891"
$ws.Range("M14").Value = "[43.288461592940564, 78.74098301281725]"
$ws.Range("N14").Value = [double]"1.282549222381135e-08"
$ws.Range("O14").Value = [double]"1.282549222381135e-08"
$ws.Range("P14").Value = [double]"-2.943474197958005"
$ws.Range("Q14").Value = "[-3.257947937056083, -2.6290004588599265]"
$ws.Range("T14").Value = [double]"54.08844794733661"
$ws.Range("U14").Value = "[43.73042204410327, 64.44647385056996]"
$ws.Range("V14").Value = [double]"1.045830089196897e-13"
$ws.Range("W14").Value = [double]"1.045830089196897e-13"
$ws.Range("X14").Value = [double]"11.94126126126152"
$ws.Range("Y14").Value = [double]"10.66548548548571"
$ws.Range("Z14").Value = [double]"13.21703703703733"

Write-Output "done"
